$d = $word.ActiveDocument

# The "STOP - DISCUSS..." quote paragraph ends with a run whose text is
# " BREAK DOWN THE PROBLEM". Append a new run containing just "?"
# immediately after it (within the same paragraph).
$rng = $d.Content
$rng.Find.Execute("THE VIGENERECIPHER CLASS TO BREAK DOWN THE PROBLEM")
$rng.Collapse(0)
$rng.InsertAfter("?")
